$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "FailureResponse" sample XML payload that was sitting in D2
# (the underlying shared string is dropped once nothing references it).
$ws.Range("D2").Value = ""

# Move the saved cursor/selection from D13 to A13.
$ws.Range("A13").Select()
